# Insert two new rows at the top of the "Camote" (Zapallo) price block
# (rows 676-677 in the original sheet), pushing the existing rows 676-737
# down to 678-739, then populate the two new rows with the latest weekly
# price report for "Provincia de Melipilla".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before the current row 676; everything from the
# old row 676 onward shifts down by two rows (old 676 -> 678, ... old 737 -> 739).
$ws.Range("A676:A677").EntireRow.Insert()

# New row 676: "1a (guarda)" quality, origin "Provincia de Melipilla"
$ws.Range("A676").Value = 8
$ws.Range("B676").Value = "Terminal La Palmera de La Serena"
$ws.Range("C676").Value = "Coquimbo"
$ws.Range("D676").Value = 44769
$ws.Range("E676").Value = 4
$ws.Range("F676").Value = 100112045
$ws.Range("G676").Value = "Zapallo"
$ws.Range("H676").Value = "Camote"
$ws.Range("I676").Value = "1a (guarda)"
$ws.Range("J676").Value = 1800
$ws.Range("K676").Value = 1200
$ws.Range("L676").Value = 1300
$ws.Range("M676").Value = 1250
$ws.Range("N676").Value = "$/kilo (volumen en unidades)"
$ws.Range("O676").Value = "Provincia de Melipilla"
$ws.Range("P676").Value = 1250
$ws.Range("Q676").Value = 1
$ws.Range("R676").Value = "Hortaliza"

# New row 677: "2a (guarda)" quality, origin "Provincia de Melipilla"
$ws.Range("A677").Value = 8
$ws.Range("B677").Value = "Terminal La Palmera de La Serena"
$ws.Range("C677").Value = "Coquimbo"
$ws.Range("D677").Value = 44769
$ws.Range("E677").Value = 4
$ws.Range("F677").Value = 100112045
$ws.Range("G677").Value = "Zapallo"
$ws.Range("H677").Value = "Camote"
$ws.Range("I677").Value = "2a (guarda)"
$ws.Range("J677").Value = 1000
$ws.Range("K677").Value = 1100
$ws.Range("L677").Value = 1150
$ws.Range("M677").Value = 1125
$ws.Range("N677").Value = "$/kilo (volumen en unidades)"
$ws.Range("O677").Value = "Provincia de Melipilla"
$ws.Range("P677").Value = 1125
$ws.Range("Q677").Value = 1
$ws.Range("R677").Value = "Hortaliza"
